$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Test Cases": swap D8 (FAIL->SKIP) and D14 (SKIP->FAIL)
# ---------------------------------------------------------------
$wsTestCases = $wb.Worksheets.Item("Test Cases")
$wsTestCases.Range("D8").Value = "SKIP"
$wsTestCases.Range("D14").Value = "FAIL"

# ---------------------------------------------------------------
# Sheet "TestCase_A7": D1 Results->PASS, D4 FAIL->PASS, D5 FAIL->PASS
# ---------------------------------------------------------------
$wsA7 = $wb.Worksheets.Item("TestCase_A7")
$wsA7.Range("D1").Value = "PASS"
$wsA7.Range("D4").Value = "PASS"
$wsA7.Range("D5").Value = "PASS"

# ---------------------------------------------------------------
# Sheet "TestCase_A8": swap casing of Transaction@2/transaction@2
# values in column B, and set column D results to PASS
# ---------------------------------------------------------------
$wsA8 = $wb.Worksheets.Item("TestCase_A8")
$wsA8.Range("B2").Value = "Transaction@2"
$wsA8.Range("D2").Value = "PASS"
$wsA8.Range("B3").Value = "transaction@2"
$wsA8.Range("D3").Value = "PASS"
$wsA8.Range("B4").Value = "transaction@2"
$wsA8.Range("D4").Value = "PASS"
$wsA8.Range("D6").Value = "PASS"
$wsA8.Range("B7").Value = "Transaction@2"
$wsA8.Range("D7").Value = "PASS"

# Update the saved selection/view for TestCase_A8 (D1 -> D2)
$wsA8.Activate()
$wsA8.Range("D2").Select()

# ---------------------------------------------------------------
# Update the saved selection/view for the "Test Cases" sheet:
# remove the scrolled topLeftCell and move the selection to C2:C18.
# Activate this sheet last so it remains the selected tab.
# ---------------------------------------------------------------
$wsTestCases.Activate()
$wsTestCases.Range("C2:C18").Select()
